# Femacal de La Calera - Berenjena: insert a new weekly record as row 242,
# pushing the existing rows 242-307 down to 243-308.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 242. Excel will shift rows
# 242:307 down to 243:308 and the new blank row 242 inherits formatting
# (including the date style on column D) from the row above it.
$ws.Rows.Item(242).Insert()

# Row 242 now sits between the (unchanged) row 241 and what used to be
# row 242 (now row 243). Populate the constant columns from the row that
# follows it (row 243), which still carries the original record's fixed
# attributes (market, region, product, quality, unit, origin, etc.).
$ws.Range("A242").Value = 3
$ws.Range("B242").Value = "Femacal de La Calera"
$ws.Range("C242").Value = "Coquimbo"
$ws.Range("D242").Value = 44736
$ws.Range("E242").Value = 5
$ws.Range("F242").Value = 100112001
$ws.Range("G242").Value = "Berenjena"
$ws.Range("H242").Value = "Sin especificar"
$ws.Range("I242").Value = "Primera"
$ws.Range("J242").Value = 130
$ws.Range("K242").Value = 7000
$ws.Range("L242").Value = 7500
$ws.Range("M242").Value = 7250
$ws.Range("N242").Value = "$/caja 60 unidades"
$ws.Range("O242").Value = "Región de Arica y Parinacota"
$ws.Range("P242").Value = 121
$ws.Range("Q242").Value = 60
$ws.Range("R242").Value = "Hortaliza"
